# Daily attendance processing - 2025-11-06 07:22:20
#
# - Reorders the "Recorded By" (column G) email lists for several sessions
#   (same set of recorders, new display order from the sync job).
# - Refreshes the Year-2 summary counters (recorded/missing session counts
#   and coverage / average-attendance percentages).
# - Marks the Year 2 / A3 / PHARMACOLOGY session #2 (row 46) as processed:
#   it flips from "Not Recorded" (pink) to "Recorded" (green) with the
#   recorder and the attendance tally filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder "Recorded By" (column G) email lists -------------------------
$gValues = @{
    2 = 'servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
    3 = 'servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
    7 = 'Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
    18 = 'servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
    19 = 'servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
    23 = 'Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
    34 = 'servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
    35 = 'servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, System, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
    50 = 'servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
    51 = 'servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, System, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
    66 = 'servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
    67 = 'servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
    77 = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
    82 = 'servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
    83 = 'servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
    93 = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
    98 = 'servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
    103 = 'Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
    114 = 'servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg'
    119 = 'Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
}
foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

# --- Refresh summary counters ----------------------------------------------
# Recorded / Missing session counters near the top of the sheet
$ws.Range("L6").Value = 27   # Recorded Sessions (was 26)
$ws.Range("L7").Value = 12   # Missing Sessions (was 13)

# Coverage % / Average Attendance % are stored as plain text, not numbers,
# so a formula-literal is entered and then collapsed to a literal value via
# a Copy / PasteSpecial(values) round-trip. This avoids Excel's automatic
# "look like a percentage" -> numeric-percent conversion, which would change
# the cell's storage type and formatting away from the plain text used here.
function Set-TextValue($rangeAddress, $text) {
    $target = $ws.Range($rangeAddress)
    $target.Formula = '="' + $text + '"'
    $target.Copy()
    $target.PasteSpecial(-4163)   # xlPasteValues
}

Set-TextValue "L9" "21.1%"    # Coverage % (was 20.3%)
Set-TextValue "L10" "56.2%"   # Average Attendance % (was 58.1%)

# Per-subject breakdown row for Year 2 / A3 (row 17)
$ws.Range("O17").Value = 3    # was 2
$ws.Range("P17").Value = 2    # was 3
Set-TextValue "R17" "18.8%"   # was 12.5%
Set-TextValue "S17" "50.2%"   # was 72.5%

# --- Row 46: Year 2 / A3 / PHARMACOLOGY #2 is now recorded ------------------
# Copy the green "Recorded" formatting from row 2 (format-only paste, so the
# existing A46:I46 values are left untouched), then fill in the recorder and
# the updated attendance counts/status.
$ws.Range("A2:I2").Copy()
$ws.Range("A46:I46").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G46").Value = "nourhan.mostafa@med.asu.edu.eg"
$ws.Range("H46").Value = "12/220"
$ws.Range("I46").Value = "Recorded"
